# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary rows (10,11,12): update scores ---
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 44
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "42/112"

# Apply the mtitleStyle (same style used by the header row 9 / "No." "Marking" "Total" labels)
# to A10, A11, A12 without disturbing their text, by copying format from A9.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# --- Remove the third answer block (columns G/H) entirely, and trim the
# second answer block (columns D/E) so only the header row (15) and the
# first two answer rows (16,17,18) still show a "Student Ans"/"Correct Ans" pair ---
$ws.Range("F5:H40").Clear()
$ws.Range("D19:E40").Clear()

# --- Fill in the student's answers for the first answer block (column A),
# marking each as correct or incorrect depending on whether it matches the
# correct answer already stored in column B. correctStyle is borrowed from
# B10, incorrectStyle from C10 (so the existing cell style entries are reused). ---
$ws.Range("B10").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Option A"

$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Option B"

$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Option C"

$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Option A"

$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = "Option B"

$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A33").Value = "Option D"

$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = "Option A"

# Student's answer for the remaining (second) answer block (column D)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Option D"

# Incorrect answers
$ws.Range("C10").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = "Option D"

$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Option D"
